$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5605909413
$ws.Range("C2").Value = -224.73419782
$ws.Range("D2").Value = -225.29478876
$ws.Range("E2").Value = -224.5092728215

$ws.Range("B3").Value = -0.5691981028999999
$ws.Range("C3").Value = -224.66645069
$ws.Range("D3").Value = -225.23564879
$ws.Range("E3").Value = -224.5092728215

$ws.Range("B4").Value = -0.5711951382
$ws.Range("C4").Value = -224.64262073
$ws.Range("D4").Value = -225.21381587
$ws.Range("E4").Value = -224.5092728215
